{"js": "// Append a new block of paragraphs (a blank-line gap, an explanatory\n// sentence, another blank line, and a third-wave R code chunk) after the\n// last paragraph of the document body - \"Comments from Ian and Nicky\n// incorporated\".\nconst newParagraphs = [\n  \"\",\n  \"\",\n  \"\",\n  \"\",\n  \"Below we perform a third wave of the process. Everything is as done for the second wave, apart the fact that we use $100$ repetitions here, instead of $50$.\",\n  \"\",\n  \"```{r}\",\n  \"new_new_results <- list()\",\n  \"with_progress({\",\n  \"  p <- progressor(nrow(initial_points))\",\n  \"for (i in 1:nrow(new_new_points)) {\",\n  \"  model_out <- get_results(unlist(new_new_points[i,]), nreps = 100, outs = c(\\\"I\\\", \\\"R\\\"), \",\n  \"                           times = c(25, 40, 100, 200))\",\n  \"  new_new_results[[i]] <- model_out\",\n  \"  p(message = sprintf(\\\"Run %g\\\", i))\",\n  \"}\",\n  \"})\",\n  \"wave2 <- data.frame(do.call('rbind', new_new_results))\",\n  \"new_new_all_training <- wave2[1:10000,]\",\n  \"new_new_all_valid <- wave2[10001:15000,]\",\n  \"new_new_stoch_emulators <- variance_emulator_from_data(new_new_all_training, output_names, ranges, \",\n  \"                                                  check.ranges=TRUE)\",\n  \"new_new_new_points <- generate_new_runs(c(new_new_stoch_emulators, \",\n  \"                                          new_stoch_emulators, stoch_emulators), 150, targets)\",\n  \"```\",\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nlet anchor = paragraphs.items[paragraphs.items.length - 1];\nfor (const text of newParagraphs) {\n  anchor = anchor.insertParagraph(text, \"After\");\n}\nawait context.sync();\n", "ps1": "# Append a new block of paragraphs (a blank-line gap, an explanatory\n# sentence, another blank line, and a third-wave R code chunk) after the\n# last paragraph of the document body - \"Comments from Ian and Nicky\n# incorporated\".\n$doc = $word.ActiveDocument\n\n$newParagraphs = @(\n    '',\n    '',\n    '',\n    '',\n    'Below we perform a third wave of the process. Everything is as done for the second wave, apart the fact that we use $100$ repetitions here, instead of $50$.',\n    '',\n    '```{r}',\n    'new_new_results <- list()',\n    'with_progress({',\n    '  p <- progressor(nrow(initial_points))',\n    'for (i in 1:nrow(new_new_points)) {',\n    '  model_out <- get_results(unlist(new_new_points[i,]), nreps = 100, outs = c(\"I\", \"R\"), ',\n    '                           times = c(25, 40, 100, 200))',\n    '  new_new_results[[i]] <- model_out',\n    '  p(message = sprintf(\"Run %g\", i))',\n    '}',\n    '})',\n    'wave2 <- data.frame(do.call(''rbind'', new_new_results))',\n    'new_new_all_training <- wave2[1:10000,]',\n    'new_new_all_valid <- wave2[10001:15000,]',\n    'new_new_stoch_emulators <- variance_emulator_from_data(new_new_all_training, output_names, ranges, ',\n    '                                                  check.ranges=TRUE)',\n    'new_new_new_points <- generate_new_runs(c(new_new_stoch_emulators, ',\n    '                                          new_stoch_emulators, stoch_emulators), 150, targets)',\n    '```'\n)\n\nforeach ($t in $newParagraphs) {\n    $rng = $doc.Paragraphs.Last.Range\n    $rng.InsertParagraphAfter()\n    if ($t -ne '') {\n        $doc.Paragraphs.Last.Range.InsertAfter($t)\n    }\n}\n"}
